$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Price" column keeps being stored as text (it holds values
# like "64.215.70" / "0.740" / "1.00" that must not be coerced to numbers
# and lose their original formatting/precision).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.215.70"
$ws.Range("E2").Value = "  +0.26%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.116.81"
$ws.Range("E3").Value = "  -0.49%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.12%  "

# Row 5 - BNB
$ws.Range("D5").Value = "591.06"
$ws.Range("E5").Value = "  +0.62%  "

# Row 6 - Solana
$ws.Range("D6").Value = "154.02"
$ws.Range("E6").Value = "  +5.42%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.118.55"
$ws.Range("E8").Value = "  -0.33%  "

# Row 9 - XRP
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  +0.87%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.27%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "5.98"
$ws.Range("E11").Value = "  +2.65%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +0.93%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -0.14%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "38.04"
$ws.Range("E14").Value = "  +2.36%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.630.80"
$ws.Range("E15").Value = "  -0.91%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  -1.62%  "

# Row 17 - Polkadot
$ws.Range("D17").Value = "7.22"
$ws.Range("E17").Value = "  +2.11%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "64.032.51"
$ws.Range("E18").Value = "  +0.27%  "

# Row 19 - WrappedEther
$ws.Range("D19").Value = "3.117.53"
$ws.Range("E19").Value = "  -1.10%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "472.77"
$ws.Range("E20").Value = "  +2.09%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "14.84"
$ws.Range("E21").Value = "  +4.01%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "0.740"
$ws.Range("E22").Value = "  +1.76%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "7.61"
$ws.Range("E23").Value = "  +3.17%  "

# Row 24 - Fetch.AI
$ws.Range("E24").Value = "  +5.30%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("D25").Value = "13.30"

# Row 26 - Litecoin
$ws.Range("D26").Value = "81.87"
$ws.Range("E26").Value = "  +1.46%  "

# Row 27 - was Dai, now RenderToken
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "9.99"
$ws.Range("E27").Value = "  +5.81%  "

# Row 28 - was RenderToken, now Dai
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.21%  "

# Row 29 - NEARProtocol
$ws.Range("D29").Value = "7.45"
$ws.Range("E29").Value = "  +4.78%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +1.57%  "

# Row 31 - was ImmutableX, now FirstDigitalUSD
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.09%  "

# Row 32 - was FirstDigitalUSD, now ImmutableX
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "2.21"
$ws.Range("E32").Value = "  +1.08%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.117"
$ws.Range("E33").Value = "  +7.42%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "27.61"
$ws.Range("E34").Value = "  +2.59%  "

# Row 35 - PEPE
$ws.Range("D35").Value = "0.0₃0863"
$ws.Range("E35").Value = "  +1.32%  "

# Row 36 - Mantle
$ws.Range("E36").Value = "  +2.04%  "

# Row 37 - dogwifhat
$ws.Range("D37").Value = "3.45"
$ws.Range("E37").Value = "  +4.41%  "

# Row 38 - Filecoin
$ws.Range("D38").Value = "6.18"
$ws.Range("E38").Value = "  +2.70%  "

# Row 39 - Stacks
$ws.Range("D39").Value = "2.27"
$ws.Range("E39").Value = "  -1.42%  "

# Row 40 - was Cosmos, now Bittensor
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "458.16"
$ws.Range("E40").Value = "  +4.69%  "

# Row 41 - was Bittensor, now Cosmos
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "9.34"
$ws.Range("E41").Value = "  +5.37%  "

# Row 42 - OKB
$ws.Range("E42").Value = "  -0.52%  "

# Row 43 - TheGraph
$ws.Range("E43").Value = "  +2.14%  "

# Row 44 - VeChain
$ws.Range("D44").Value = "0.0372"
$ws.Range("E44").Value = "  +0.44%  "

# Row 45 - Maker
$ws.Range("D45").Value = "2.864.05"
$ws.Range("E45").Value = "  -1.36%  "

# Row 46 - Kaspa
$ws.Range("D46").Value = "0.111"
$ws.Range("E46").Value = "  +3.05%  "

# Row 47 - Arweave
$ws.Range("D47").Value = "39.75"
$ws.Range("E47").Value = "  +0.39%  "

# Row 48 - Monero
$ws.Range("D48").Value = "131.02"
$ws.Range("E48").Value = "  +3.56%  "

# Row 49 - InjectiveProtocol
$ws.Range("D49").Value = "25.63"
$ws.Range("E49").Value = "  +6.27%  "

# Row 50 - was ThetaToken, now USDe
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  +0.03%  "

# Row 51 - was USDe, now ThetaToken
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "2.28"
$ws.Range("E51").Value = "  +4.46%  "
